# Scheduled runner update: refresh Leve profit calculations with latest market prices
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 2946.6155
$ws.Range("I38").Value = 992.3333
$ws.Range("K38").Value = 2976.9999
$ws.Range("M38").Value = -2604.9999
$ws.Range("H39").Value = 1004
$ws.Range("J39").Value = 2000
$ws.Range("L39").Value = 6000
$ws.Range("N39").Value = -6592
$ws.Range("H42").Value = 228.25
$ws.Range("J42").Value = 380
$ws.Range("L42").Value = 1140
$ws.Range("N42").Value = -1600
$ws.Range("H51").Value = 4957.875
$ws.Range("J51").Value = 4957.875
$ws.Range("L51").Value = 4957.875
$ws.Range("N51").Value = -5925.875
$ws.Range("H62").Value = 6626.5
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 6626.5
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H70").Value = 3450
$ws.Range("J70").Value = 3450
$ws.Range("L70").Value = 10350
$ws.Range("N70").Value = -10890
$ws.Range("H73").Value = 3450
$ws.Range("J73").Value = 3450
$ws.Range("L73").Value = 10350
$ws.Range("N73").Value = -12222
$ws.Range("H112").Value = 5816035.5
$ws.Range("J112").Value = 6026121
$ws.Range("L112").Value = 18078363
$ws.Range("N112").Value = -18080579
$ws.Range("H116").Value = 2839.2222
$ws.Range("J116").Value = 2383.3333
$ws.Range("L116").Value = 2383.3333
$ws.Range("N116").Value = -9267.3333
$ws.Range("H118").Value = 1127.1765
$ws.Range("I118").Value = 1067.0769
$ws.Range("K118").Value = 3201.2307
$ws.Range("M118").Value = -1544.2307
$ws.Range("H132").Value = 5045.6
$ws.Range("I132").Value = 5266.4443
$ws.Range("K132").Value = 15799.3329
$ws.Range("M132").Value = -13269.3329
$ws.Range("H138").Value = 7409878
$ws.Range("I138").Value = 1338.5714
$ws.Range("J138").Value = 10755670
$ws.Range("K138").Value = 4015.7142
$ws.Range("L138").Value = 32267010
$ws.Range("M138").Value = 1124.2858
$ws.Range("N138").Value = -32277290

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 148.75
$ws.Range("I4").Value = 72.5
$ws.Range("J4").Value = 225
$ws.Range("K4").Value = 72.5
$ws.Range("L4").Value = 225
$ws.Range("M4").Value = 43.5
$ws.Range("N4").Value = -457
$ws.Range("H39").Value = 9631.666999999999
$ws.Range("I39").Value = 9631.666999999999
$ws.Range("K39").Value = 9631.666999999999
$ws.Range("M39").Value = -9111.666999999999
$ws.Range("H97").Value = 851.46155
$ws.Range("I97").Value = 961.76
$ws.Range("J97").Value = 654.5
$ws.Range("K97").Value = 961.76
$ws.Range("L97").Value = 654.5
$ws.Range("M97").Value = -465.76
$ws.Range("N97").Value = -1646.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1613.591
$ws.Range("I20").Value = 1565.4375
$ws.Range("J20").Value = 1742
$ws.Range("K20").Value = 1565.4375
$ws.Range("L20").Value = 1742
$ws.Range("M20").Value = -1318.4375
$ws.Range("N20").Value = -2236
$ws.Range("H134").Value = 2764.5
$ws.Range("I134").Value = 2414.5
$ws.Range("J134").Value = 3377
$ws.Range("K134").Value = 7243.5
$ws.Range("L134").Value = 10131
$ws.Range("M134").Value = -4708.5
$ws.Range("N134").Value = -15201

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 13485.429
$ws.Range("I99").Value = 13485.429
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 13485.429
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -11987.429
$ws.Range("N99").ClearContents()
$ws.Range("H105").Value = 16811.715
$ws.Range("I105").Value = 1466.6666
$ws.Range("J105").Value = 28320.5
$ws.Range("K105").Value = 1466.6666
$ws.Range("L105").Value = 28320.5
$ws.Range("M105").Value = 280.3334
$ws.Range("N105").Value = -31814.5
$ws.Range("H126").Value = 13485.429
$ws.Range("I126").Value = 13485.429
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 40456.287
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -37986.287
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 227326.11
$ws.Range("J132").Value = 6136.143
$ws.Range("L132").Value = 18408.429
$ws.Range("N132").Value = -23468.429
$ws.Range("H134").Value = 2581.6
$ws.Range("I134").Value = 2581.6
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 7744.799999999999
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -5209.799999999999
$ws.Range("N134").ClearContents()
$ws.Range("H141").Value = 187852.14
$ws.Range("J141").Value = 187852.14
$ws.Range("L141").Value = 187852.14
$ws.Range("N141").Value = -198212.14

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 83.5
$ws.Range("J38").Value = 85.2
$ws.Range("L38").Value = 255.6
$ws.Range("N38").Value = -949.6
$ws.Range("H44").Value = 1824.8889
$ws.Range("I44").Value = 1740.5
$ws.Range("K44").Value = 5221.5
$ws.Range("M44").Value = -4823.5
$ws.Range("H131").Value = 1665.1316
$ws.Range("J131").Value = 1654.7297
$ws.Range("L131").Value = 4964.189100000001
$ws.Range("N131").Value = -15044.1891
$ws.Range("H133").Value = 3160.6
$ws.Range("I133").Value = 3160.6
$ws.Range("K133").Value = 9481.799999999999
$ws.Range("M133").Value = -4421.799999999999
$ws.Range("H134").Value = 5720.92
$ws.Range("I134").Value = 1472.0588
$ws.Range("K134").Value = 4416.1764
$ws.Range("M134").Value = 653.8235999999997
$ws.Range("H137").Value = 3384.8125
$ws.Range("J137").Value = 3419.8462
$ws.Range("L137").Value = 10259.5386
$ws.Range("N137").Value = -20459.5386

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 30099
$ws.Range("J49").Value = 30747.5
$ws.Range("L49").Value = 30747.5
$ws.Range("N49").Value = -31115.5
$ws.Range("H80").Value = 3047.7334
$ws.Range("I80").Value = 2996.25
$ws.Range("K80").Value = 2996.25
$ws.Range("M80").Value = -1998.25
$ws.Range("H83").Value = 3047.7334
$ws.Range("I83").Value = 2996.25
$ws.Range("K83").Value = 14981.25
$ws.Range("M83").Value = -9989.25

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 25000
$ws.Range("I42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("M42").ClearContents()
$ws.Range("H49").Value = 25000
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").ClearContents()
$ws.Range("H68").Value = 3389
$ws.Range("J68").Value = 5625
$ws.Range("L68").Value = 5625
$ws.Range("N68").Value = -7123
$ws.Range("H71").Value = 3389
$ws.Range("J71").Value = 5625
$ws.Range("L71").Value = 28125
$ws.Range("N71").Value = -35613

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H38").Value = 5020000
$ws.Range("I38").Value = 10007500
$ws.Range("J38").Value = 32499.5
$ws.Range("K38").Value = 10007500
$ws.Range("L38").Value = 32499.5
$ws.Range("M38").Value = -10007027
$ws.Range("N38").Value = -33445.5
$ws.Range("H42").Value = 50000
$ws.Range("J42").Value = 50000
$ws.Range("L42").Value = 50000
$ws.Range("N42").Value = -50756
$ws.Range("H100").Value = 76924250
$ws.Range("I100").Value = 142858820
$ws.Range("K100").Value = 285717640
$ws.Range("M100").Value = -285717099
$ws.Range("H122").Value = 41713.08
$ws.Range("I122").Value = 51333.8
$ws.Range("J122").Value = 3230.2
$ws.Range("K122").Value = 154001.4
$ws.Range("L122").Value = 9690.599999999999
$ws.Range("M122").Value = -151551.4
$ws.Range("N122").Value = -14590.6
$ws.Range("H126").Value = 5175.467
$ws.Range("I126").Value = 6259.1665
$ws.Range("J126").Value = 840.6667
$ws.Range("K126").Value = 18777.4995
$ws.Range("L126").Value = 2522.0001
$ws.Range("M126").Value = -16307.4995
$ws.Range("N126").Value = -7462.0001

